$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualiza dados em 28-12-2017
$ws.Range("B2").Value = 43097.26338454861
$ws.Range("B3").Value = 43097.27128064814
$ws.Range("B4").Value = 43097.27089083333
$ws.Range("B5").Value = 43097.30682739583
